$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRange, $text)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "64.987.65"
Set-TextCell $ws.Range("E2") "  +0.97%  "

Set-TextCell $ws.Range("D3") "3.374.05"
Set-TextCell $ws.Range("E3") "  +0.58%  "

Set-TextCell $ws.Range("D4") "1.00"
Set-TextCell $ws.Range("E4") "  +0.06%  "

Set-TextCell $ws.Range("D5") "554.08"
Set-TextCell $ws.Range("E5") "  -0.30%  "

Set-TextCell $ws.Range("D6") "174.03"
Set-TextCell $ws.Range("E6") "  -0.77%  "

Set-TextCell $ws.Range("D7") "0.631"
Set-TextCell $ws.Range("E7") "  +2.03%  "

Set-TextCell $ws.Range("D8") "3.365.05"
Set-TextCell $ws.Range("E8") "  +0.80%  "

Set-TextCell $ws.Range("D9") "1.00"
Set-TextCell $ws.Range("E9") "  +0.02%  "

Set-TextCell $ws.Range("D10") "0.173"
Set-TextCell $ws.Range("E10") "  +5.37%  "

Set-TextCell $ws.Range("D11") "0.636"
Set-TextCell $ws.Range("E11") "  +1.25%  "

Set-TextCell $ws.Range("D12") "53.50"
Set-TextCell $ws.Range("E12") "  -1.90%  "

Set-TextCell $ws.Range("D13") "0.0000278"
Set-TextCell $ws.Range("E13") "  +1.74%  "

Set-TextCell $ws.Range("D14") "9.15"
Set-TextCell $ws.Range("E14") "  +0.62%  "

Set-TextCell $ws.Range("D15") "3.917.48"
Set-TextCell $ws.Range("E15") "  +0.62%  "

Set-TextCell $ws.Range("D16") "18.31"
Set-TextCell $ws.Range("E16") "  -0.38%  "

Set-TextCell $ws.Range("D17") "3.391.70"
Set-TextCell $ws.Range("E17") "  +1.17%  "

Set-TextCell $ws.Range("D18") "0.118"
Set-TextCell $ws.Range("E18") "  -0.04%  "

Set-TextCell $ws.Range("D19") "65.053.68"
Set-TextCell $ws.Range("E19") "  +1.16%  "

Set-TextCell $ws.Range("D20") "11.83"
Set-TextCell $ws.Range("E20") "  -0.11%  "

Set-TextCell $ws.Range("D21") "0.996"
Set-TextCell $ws.Range("E21") "  +1.25%  "

Set-TextCell $ws.Range("D22") "455.09"
Set-TextCell $ws.Range("E22") "  -0.18%  "

Set-TextCell $ws.Range("D23") "4.86"
Set-TextCell $ws.Range("E23") "  +0.00%  "

Set-TextCell $ws.Range("D24") "14.16"
Set-TextCell $ws.Range("E24") "  +6.13%  "

Set-TextCell $ws.Range("D25") "4.07"
Set-TextCell $ws.Range("E25") "  -0.29%  "

Set-TextCell $ws.Range("D26") "87.37"
Set-TextCell $ws.Range("E26") "  +1.77%  "

Set-TextCell $ws.Range("D27") "2.87"
Set-TextCell $ws.Range("E27") "  +0.86%  "

Set-TextCell $ws.Range("D28") "10.68"
Set-TextCell $ws.Range("E28") "  -2.84%  "

Set-TextCell $ws.Range("D29") "8.69"
Set-TextCell $ws.Range("E29") "  -1.00%  "

Set-TextCell $ws.Range("D30") "31.05"
Set-TextCell $ws.Range("E30") "  +3.48%  "

Set-TextCell $ws.Range("D31") "6.51"
Set-TextCell $ws.Range("E31") "  -2.41%  "

Set-TextCell $ws.Range("D32") "63.29"
Set-TextCell $ws.Range("E32") "  +7.72%  "

Set-TextCell $ws.Range("D33") "11.45"
Set-TextCell $ws.Range("E33") "  -0.29%  "

Set-TextCell $ws.Range("D34") "576.69"
Set-TextCell $ws.Range("E34") "  -1.49%  "

Set-TextCell $ws.Range("D35") "0.107"
Set-TextCell $ws.Range("E35") "  -0.79%  "

Set-TextCell $ws.Range("D36") "1.00"
Set-TextCell $ws.Range("E36") "  -0.02%  "

Set-TextCell $ws.Range("D37") "3.59"
Set-TextCell $ws.Range("E37") "  +2.17%  "

Set-TextCell $ws.Range("D38") "0.143"
Set-TextCell $ws.Range("E38") "  +1.60%  "

Set-TextCell $ws.Range("D39") "35.61"
Set-TextCell $ws.Range("E39") "  -0.38%  "

Set-TextCell $ws.Range("D40") "0.372"
Set-TextCell $ws.Range("E40") "  -0.81%  "

Set-TextCell $ws.Range("D41") "0.0₃0735"
Set-TextCell $ws.Range("E41") "  -2.74%  "

Set-TextCell $ws.Range("D42") "3.096.65"
Set-TextCell $ws.Range("E42") "  +0.07%  "

Set-TextCell $ws.Range("D43") "0.0415"
Set-TextCell $ws.Range("E43") "  +1.16%  "

Set-TextCell $ws.Range("D44") "2.74"
Set-TextCell $ws.Range("E44") "  -1.63%  "

Set-TextCell $ws.Range("B45") "ApeXProtocol"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws.Range("D45") "3.19"
Set-TextCell $ws.Range("E45") "  +0.32%  "

Set-TextCell $ws.Range("B46") "Stellar"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D46") "0.134"
Set-TextCell $ws.Range("E46") "  +2.34%  "

Set-TextCell $ws.Range("B47") "Fetch.AI"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Range("D47") "2.44"
Set-TextCell $ws.Range("E47") "  -4.00%  "

Set-TextCell $ws.Range("D48") "1.00"
Set-TextCell $ws.Range("E48") "  +0.13%  "

Set-TextCell $ws.Range("D49") "140.80"
Set-TextCell $ws.Range("E49") "  +4.02%  "

Set-TextCell $ws.Range("D50") "2.52"
Set-TextCell $ws.Range("E50") "  -2.16%  "

Set-TextCell $ws.Range("D51") "8.30"
Set-TextCell $ws.Range("E51") "  -0.71%  "
